# Applies the Google-Sheet-driven content refresh:
#   - fills in the previously-empty I32/J32 (media_hero/media_alt) cells
#   - appends two new "songs" collection rows (145, 146) with their full data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: Virus tagline gains a media_hero link + media_alt text ---
$ws.Range("I32").Value2 = "https://drive.google.com/file/d/1I6WY4dTmIMX3YbyvpTjuBBFBJ53twXxt/view?usp=drive_link"
$ws.Range("J32").Value2 = "Virus Database has been outdated"

# --- Row 145: new "songs" entry, "Mind the Gap" ---
$ws.Range("A145").Value2 = "songs"
$ws.Range("B145").Value2 = "Mind the Gap"
$ws.Range("C145").Value2 = "distance-step-gap-return-deadend-consistency-trust"
$ws.Range("D145").Value2 = "songs/distance-step-gap-return-deadend-consistency-trust/"
$ws.Range("E145").Value2 = "MusicComposition"
$ws.Range("F145").Value2 = ""
$ws.Range("G145").Value2 = ""
$ws.Range("H145").Value2 = "mind the gap, watch your step, keep your distance"
$ws.Range("I145").Value2 = "https://drive.google.com/file/d/1nDd-NpjMKQ5TT1ARBL9hQEuqHhGCRgnk/view?usp=drive_link"
$ws.Range("J145").Value2 = "Mind the Gap, Soul Feedback"
$ws.Range("K145").Value2 = ""
$ws.Range("L145").Value2 = "Counterfeit"
$ws.Range("M145").Value2 = ""
$ws.Range("N145").Value2 = @"
GAP
I don't believe in words anymore,
I mirror actions and consistency,
and since I don't trust you at all
You better stay in line:
You better keep your distance mister,
You better watch your step,
You better mind the gap.
Our pictures
Will disappear
As if we were never together
In the first place.
So play the role
Of a good
Stranger:
For you
there's just One way,
No return,
Dead end.
"@

# --- Row 146: new "songs" entry, "Opposite" ---
$ws.Range("A146").Value2 = "songs"
$ws.Range("B146").Value2 = "Opposite"
$ws.Range("C146").Value2 = "opposite-directions-train-thoughts-obliviousness-avatar"
$ws.Range("D146").Value2 = "songs/opposite-directions-train-thoughts-obliviousness-avatar/"
$ws.Range("E146").Value2 = "MusicComposition"
$ws.Range("F146").Value2 = ""
$ws.Range("G146").Value2 = ""
$ws.Range("H146").Value2 = "Opposite Directions"
$ws.Range("I146").Value2 = "https://drive.google.com/file/d/1y-oTWR0YiQwTstdUHOo3Jz40_IAlkNWl/view?usp=drive_link"
$ws.Range("J146").Value2 = "Opposite Directions, Soul Feedback"
$ws.Range("K146").Value2 = ""
$ws.Range("L146").Value2 = "Counterfeit"
$ws.Range("M146").Value2 = ""
$ws.Range("N146").Value2 = @"
OPPOSITE DIRECTIONS
Is not nice, not to see you, despite you being in front of me,
Poker face, black holes in your eyes,
There's no further distance than obliviousness,
Two strangers going to opposite directions,
Empty vacuum gazes and complete lack of interest,
Is when physical proximity turns damnation.
And I wasted my machine-machine-gun against your avatar...
Smashing, slashing, flashing, crashing...
Opposite Directions in the train of thoughts,
I crossed the treshold: expelled from paradise.
Opposite Directions, streams of consciousness,
Too much desperation, there's no turning back.
"@
